# Commit: "Changed GPA as number"
# The CGPA column (Y) used the shared string "--" to represent a blank/zero
# GPA for students who failed (column X = "F"). This script converts every
# such "--" text cell in column Y to the plain number 0, matching the
# cleaned-data re-export. Once no cell references the "--" shared string any
# more, the workbook writer drops it from sharedStrings.xml, which is why a
# lot of the other <v> indices shift down by one as a side effect.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 25).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 206 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 25)
    $val = $cell.Value()
    if ($val -eq "--") {
        $cell.Value = 0
    }
}

# Restore the view: scroll near the top of the sheet and select Y1, as in
# the target workbook.
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$win.ScrollRow = 1
$ws.Range("Y1").Select()
